$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 36 (Washington) updates
$ws.Range("B36").Value = 44034
$ws.Range("C36").Value = 49247
$ws.Range("D36").Value = 1468
$ws.Range("E36").Value = 1743
$ws.Range("F36").Value = 49
$ws.Range("H36").Value = 3.47
$ws.Range("K36").Value = 32066
$ws.Range("L36").Value = 1413

# Row 39 (Delaware) status message update
$ws.Range("O39").Value = "An error occurred. ... AttributeError(""'numpy.float64' object has no attribute 'split'"")"
